$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 337.22726
$ws.Range("I8").Value = 54.666668
$ws.Range("K8").Value = 164.000004
$ws.Range("M8").Value = -25.00000399999999
$ws.Range("H9").Value = 1225
$ws.Range("I9").Value = 1533.3334
$ws.Range("J9").Value = 916.6667
$ws.Range("K9").Value = 1533.3334
$ws.Range("L9").Value = 916.6667
$ws.Range("M9").Value = -1364.3334
$ws.Range("N9").Value = -1254.6667
$ws.Range("H43").Value = 522467.3
$ws.Range("I43").Value = 12970.75
$ws.Range("K43").Value = 12970.75
$ws.Range("M43").Value = -12901.75
$ws.Range("H88").Value = 20885174
$ws.Range("I88").Value = 111113430
$ws.Range("J88").Value = 63267.46
$ws.Range("K88").Value = 111113430
$ws.Range("L88").Value = 63267.46
$ws.Range("M88").Value = -111113024
$ws.Range("N88").Value = -64079.46
$ws.Range("H91").Value = 20885174
$ws.Range("I91").Value = 111113430
$ws.Range("J91").Value = 63267.46
$ws.Range("K91").Value = 111113430
$ws.Range("L91").Value = 63267.46
$ws.Range("M91").Value = -111112026
$ws.Range("N91").Value = -66075.45999999999
$ws.Range("H100").Value = 2703.6667
$ws.Range("I100").Value = 2016.5555
$ws.Range("J100").Value = 4765
$ws.Range("K100").Value = 2016.5555
$ws.Range("L100").Value = 4765
$ws.Range("M100").Value = -1475.5555
$ws.Range("N100").Value = -5847
$ws.Range("H101").Value = 1500
$ws.Range("I101").Value = 400
$ws.Range("J101").Value = 2600
$ws.Range("K101").Value = 1200
$ws.Range("L101").Value = 7800
$ws.Range("M101").Value = 422
$ws.Range("N101").Value = -11044
$ws.Range("H112").Value = 9587
$ws.Range("J112").Value = 10278.381
$ws.Range("L112").Value = 30835.143
$ws.Range("N112").Value = -33051.143
$ws.Range("H125").Value = 2922.9565
$ws.Range("I125").Value = 1788.1538
$ws.Range("J125").Value = 4398.2
$ws.Range("K125").Value = 16093.3842
$ws.Range("L125").Value = 39583.8
$ws.Range("M125").Value = -13633.3842
$ws.Range("N125").Value = -44503.8
$ws.Range("H131").Value = 1598
$ws.Range("I131").Value = 1598
$ws.Range("K131").Value = 4794
$ws.Range("M131").Value = 246
$ws.Range("H132").Value = 1986.65
$ws.Range("I132").Value = 1986.65
$ws.Range("K132").Value = 5959.950000000001
$ws.Range("M132").Value = -3429.950000000001
$ws.Range("H135").Value = 204841.89
$ws.Range("I135").Value = 257162.05
$ws.Range("J135").Value = 793.3
$ws.Range("K135").Value = 2314458.45
$ws.Range("L135").Value = 7139.7
$ws.Range("M135").Value = -2311923.45
$ws.Range("N135").Value = -12209.7
$ws.Range("H137").Value = 1437.0714
$ws.Range("I137").Value = 976.8333
$ws.Range("K137").Value = 2930.4999
$ws.Range("M137").Value = -380.4998999999998
$ws.Range("H138").Value = 4170062.2
$ws.Range("I138").Value = 2506.125
$ws.Range("J138").Value = 12505174
$ws.Range("K138").Value = 7518.375
$ws.Range("L138").Value = 37515522
$ws.Range("M138").Value = -2378.375
$ws.Range("N138").Value = -37525802
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 25000000
$ws.Range("I6").Value = 25000000
$ws.Range("K6").Value = 25000000
$ws.Range("M6").Value = -24999827
$ws.Range("H45").Value = 3846.8235
$ws.Range("I45").Value = 2440.2
$ws.Range("J45").Value = 4432.9165
$ws.Range("K45").Value = 2440.2
$ws.Range("L45").Value = 4432.9165
$ws.Range("M45").Value = -2063.2
$ws.Range("N45").Value = -5186.9165
$ws.Range("H61").Value = 9888.791999999999
$ws.Range("I61").Value = 3765.889
$ws.Range("K61").Value = 3765.889
$ws.Range("M61").Value = -3553.889
$ws.Range("H74").Value = 25042.34
$ws.Range("I74").Value = 36600.1
$ws.Range("K74").Value = 36600.1
$ws.Range("M74").Value = -35726.1
$ws.Range("H77").Value = 25042.34
$ws.Range("I77").Value = 36600.1
$ws.Range("K77").Value = 183000.5
$ws.Range("M77").Value = -178632.5
$ws.Range("H110").Value = 2414
$ws.Range("I110").Value = 2267.75
$ws.Range("K110").Value = 2267.75
$ws.Range("M110").Value = -222.75
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0
$ws.Range("H132").Value = 1256634.5
$ws.Range("I132").Value = 2178862.8
$ws.Range("K132").Value = 6536588.399999999
$ws.Range("M132").Value = -6534058.399999999
$ws.Range("H136").Value = 9888.791999999999
$ws.Range("I136").Value = 3765.889
$ws.Range("K136").Value = 11297.667
$ws.Range("M136").Value = -8747.667000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1977.2727
$ws.Range("J11").Value = 2772.6
$ws.Range("L11").Value = 2772.6
$ws.Range("N11").Value = -3052.6
$ws.Range("H88").Value = 75000
$ws.Range("J88").Value = 75000
$ws.Range("L88").Value = 75000
$ws.Range("N88").Value = -75812
$ws.Range("H91").Value = 75000
$ws.Range("J91").Value = 75000
$ws.Range("L91").Value = 75000
$ws.Range("N91").Value = -77808
$ws.Range("H105").Value = 3238.0454
$ws.Range("I105").Value = 1871.3077
$ws.Range("K105").Value = 1871.3077
$ws.Range("M105").Value = -124.3077000000001
$ws.Range("H107").Value = 90914136
$ws.Range("J107").Value = 6700
$ws.Range("L107").Value = 6700
$ws.Range("N107").Value = -10540
$ws.Range("H134").Value = 4224.851
$ws.Range("I134").Value = 1246.5454
$ws.Range("K134").Value = 3739.6362
$ws.Range("M134").Value = -1204.6362
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5680.58
$ws.Range("I31").Value = 1360.92
$ws.Range("J31").Value = 10000.24
$ws.Range("K31").Value = 1360.92
$ws.Range("L31").Value = 10000.24
$ws.Range("M31").Value = -1065.92
$ws.Range("N31").Value = -10590.24
$ws.Range("H34").Value = 5680.58
$ws.Range("I34").Value = 1360.92
$ws.Range("J34").Value = 10000.24
$ws.Range("K34").Value = 1360.92
$ws.Range("L34").Value = 10000.24
$ws.Range("M34").Value = -1158.92
$ws.Range("N34").Value = -10404.24
$ws.Range("H60").Value = 31166.5
$ws.Range("J60").Value = 31166.5
$ws.Range("L60").Value = 31166.5
$ws.Range("N60").Value = -32188.5
$ws.Range("H99").Value = 5658.36
$ws.Range("I99").Value = 3013.4285
$ws.Range("K99").Value = 3013.4285
$ws.Range("M99").Value = -1515.4285
$ws.Range("H105").Value = 5304.3
$ws.Range("I105").Value = 1019.8
$ws.Range("J105").Value = 9588.799999999999
$ws.Range("K105").Value = 1019.8
$ws.Range("L105").Value = 9588.799999999999
$ws.Range("M105").Value = 727.2
$ws.Range("N105").Value = -13082.8
$ws.Range("H126").Value = 5658.36
$ws.Range("I126").Value = 3013.4285
$ws.Range("K126").Value = 9040.2855
$ws.Range("M126").Value = -6570.2855
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5683.6665
$ws.Range("J34").Value = 10361.625
$ws.Range("L34").Value = 31084.875
$ws.Range("N34").Value = -31252.875
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 37076908
$ws.Range("I122").Value = 50052436
$ws.Range("J122").Value = 3982.2856
$ws.Range("K122").Value = 150157308
$ws.Range("L122").Value = 11946.8568
$ws.Range("M122").Value = -150154858
$ws.Range("N122").Value = -16846.8568
$ws.Range("H132").Value = 3702.5
$ws.Range("I132").Value = 2301.4443
$ws.Range("J132").Value = 5278.6875
$ws.Range("K132").Value = 6904.3329
$ws.Range("L132").Value = 15836.0625
$ws.Range("M132").Value = -4374.3329
$ws.Range("N132").Value = -20896.0625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5612.375
$ws.Range("I7").Value = 2848.5
$ws.Range("K7").Value = 2848.5
$ws.Range("M7").Value = -2736.5
$ws.Range("H82").Value = 7047252.5
$ws.Range("I82").Value = 14084507
$ws.Range("K82").Value = 14084507
$ws.Range("M82").Value = -14084146
$ws.Range("H85").Value = 7047252.5
$ws.Range("I85").Value = 14084507
$ws.Range("K85").Value = 14084507
$ws.Range("M85").Value = -14083259
$ws.Range("H100").Value = 2853.182
$ws.Range("I100").Value = 2397.889
$ws.Range("K100").Value = 2397.889
$ws.Range("M100").Value = -1856.889
$ws.Range("H126").Value = 5612.375
$ws.Range("I126").Value = 2848.5
$ws.Range("K126").Value = 8545.5
$ws.Range("M126").Value = -6075.5
$ws.Range("H133").Value = 83333.336
$ws.Range("J133").Value = 83333.336
$ws.Range("L133").Value = 83333.336
$ws.Range("N133").Value = -88393.336
$ws.Range("H136").Value = 7129.1035
$ws.Range("I136").Value = 4300.6577
$ws.Range("J136").Value = 12503.15
$ws.Range("K136").Value = 12901.9731
$ws.Range("L136").Value = 37509.45
$ws.Range("M136").Value = -10351.9731
$ws.Range("N136").Value = -42609.45
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5250
$ws.Range("J62").Value = 6500
$ws.Range("L62").Value = 6500
$ws.Range("N62").Value = -7748
$ws.Range("H65").Value = 5250
$ws.Range("J65").Value = 6500
$ws.Range("L65").Value = 32500
$ws.Range("N65").Value = -38740
$ws.Range("H126").Value = 2464.524
$ws.Range("I126").Value = 1811
$ws.Range("J126").Value = 3771.5715
$ws.Range("K126").Value = 5433
$ws.Range("L126").Value = 11314.7145
$ws.Range("M126").Value = -2963
$ws.Range("N126").Value = -16254.7145
$ws.Range("H136").Value = 26260.342
$ws.Range("I136").Value = 1268.6333
$ws.Range("J136").Value = 79814
$ws.Range("K136").Value = 3805.8999
$ws.Range("L136").Value = 239442
$ws.Range("M136").Value = -1255.8999
$ws.Range("N136").Value = -244542
